$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The rental return process had previously been scored with no points in
# column C (points achieved) for row 9 ("Rental return process"). After
# fixing the bug where items that hadn't been returned yet weren't being
# located correctly, this item now earns 5 points.
$ws.Range("C9").Value = 5

# Move the active selection to B18, matching where the grader's cursor
# ended up after making the fix.
$ws.Range("B18").Select()
